$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A1").Value = "Demo date"
$ws.Range("A2").Select()
